# Added Multiple Configuration & Util Files 2
#
# Removes the extra admin-credential rows from Sheet1, leaving just:
#   username / password
#   SuperAdmin / sadmin
#   Admin / admin123
# and moves the active selection to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old row 2 (Admin/admin) - SuperAdmin/sadmin shifts up to row 2.
$ws.Rows("2:2").Delete()

# Drop the old rows 5-7 (Admin/admin345, Adminsfds/admin345, Adminfsdfsdf/admin345),
# which after the first delete are now rows 4-6.
$ws.Rows("4:6").Delete()

# Update the selection to match the saved view state.
$ws.Range("A3").Select()
